{"js": "// Replace review text per the commit \"Added many more features\".\n// Each entry: exact original text -> new text. Using body.search with\n// matchCase so we target the precise runs described in the diff.\nconst replacements = [\n  {\n    find: \"Play Big Bucks Buffalo Gigablox for Free \\u2013 Review\",\n    replace: \"Play Big Bucks Buffalo Gigablox for Free\",\n  },\n  {\n    find: \"Gigablox feature increases potential winning combinations\",\n    replace: \"Gigablox feature increases chances of winning combinations\",\n  },\n  {\n    find: \"Atmospheric Western-style music and design\",\n    replace: \"Well-designed symbols with three-dimensional effect\",\n  },\n  {\n    find: \"Scatter symbols can trigger respins with high-value tokens\",\n    replace: \"Western-style musical theme enhances game atmosphere\",\n  },\n  {\n    find: \"Autoplay can be set for up to 1,000 spins\",\n    replace: \"Wide range of Autoplay options for convenience\",\n  },\n  {\n    find: \"Game has a lower theoretical return to player at 94%\",\n    replace: \"Theoretical return to player of 94% may be lower than some other games\",\n  },\n  {\n    find: \"Paylines are fixed at 40, limiting betting flexibility\",\n    replace: \"Limited bonus features beyond the Respin with gold token\",\n  },\n  {\n    find: \"Explore the Western desert with Big Bucks Buffalo Gigablox \\u2013 read our review, then play for free here! Gigablox symbols and Scatter Respins add excitement.\",\n    replace: \"Read our review of Big Bucks Buffalo Gigablox and play for free with this exciting slot game.\",\n  },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace review text per the commit \"Added many more features\".\n# Uses Find/Replace over the whole document Range so every matching run\n# (formatted or not) is updated while each run keeps its own formatting.\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $FindText\n    $find.Replacement.Text = $ReplaceText\n    $find.Forward = $true\n    $find.Wrap = 2          # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdReplaceAll = 2\n    $find.Execute([ref]$FindText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$ReplaceText, [ref]2) | Out-Null\n}\n\nReplace-AllText \"Play Big Bucks Buffalo Gigablox for Free \u2013 Review\" \"Play Big Bucks Buffalo Gigablox for Free\"\nReplace-AllText \"Gigablox feature increases potential winning combinations\" \"Gigablox feature increases chances of winning combinations\"\nReplace-AllText \"Atmospheric Western-style music and design\" \"Well-designed symbols with three-dimensional effect\"\nReplace-AllText \"Scatter symbols can trigger respins with high-value tokens\" \"Western-style musical theme enhances game atmosphere\"\nReplace-AllText \"Autoplay can be set for up to 1,000 spins\" \"Wide range of Autoplay options for convenience\"\nReplace-AllText \"Game has a lower theoretical return to player at 94%\" \"Theoretical return to player of 94% may be lower than some other games\"\nReplace-AllText \"Paylines are fixed at 40, limiting betting flexibility\" \"Limited bonus features beyond the Respin with gold token\"\nReplace-AllText \"Explore the Western desert with Big Bucks Buffalo Gigablox \u2013 read our review, then play for free here! Gigablox symbols and Scatter Respins add excitement.\" \"Read our review of Big Bucks Buffalo Gigablox and play for free with this exciting slot game.\"\n"}
